$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (sending cluster is now always "FAPs", recomputed TPM-based values)
# Row 2: Target cluster = ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Wnt1"
$ws.Range("C2").Value = "Ryk"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3532066666666667
$ws.Range("H2").Value = 1.05962
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 7.423863
$ws.Range("N2").Value = 22.271589
$ws.Range("O2").Value = 0.1690720838224332
$ws.Range("P2").Value = 0.1690720838224332
$ws.Range("Q2").Value = 2.62215790402
$ws.Range("R2").Value = 23.59942113618
$ws.Range("S2").Value = 0.1690720838224332
$ws.Range("T2").Value = 0.1690720838224332

# Row 3: Target cluster = FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Wnt1"
$ws.Range("C3").Value = "Ryk"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3532066666666667
$ws.Range("H3").Value = 1.05962
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 21.552384
$ws.Range("N3").Value = 64.657152
$ws.Range("O3").Value = 0.4908369772207905
$ws.Range("P3").Value = 0.4908369772207905
$ws.Range("Q3").Value = 7.61244571136
$ws.Range("R3").Value = 68.51201140223999
$ws.Range("S3").Value = 0.4908369772207905
$ws.Range("T3").Value = 0.4908369772207905

# Row 4: Target cluster = MuSCs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Wnt1"
$ws.Range("C4").Value = "Ryk"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.3532066666666667
$ws.Range("H4").Value = 1.05962
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 14.93320766666667
$ws.Range("N4").Value = 44.799623
$ws.Range("O4").Value = 0.3400909389567762
$ws.Range("P4").Value = 0.3400909389567762
$ws.Range("Q4").Value = 5.274508502584444
$ws.Range("R4").Value = 47.47057652326
$ws.Range("S4").Value = 0.3400909389567762
$ws.Range("T4").Value = 0.3400909389567762

# Remove the now-obsolete rows (old "ECs" sending-cluster rows 5-7)
$ws.Range("A5:T7").Delete()
